{"js": "const pairs = [\n  [\"2024-09-11 Wednesday\", \"2024-09-12 Thursday\"],\n  [\"61-55=6\", \"0+11=11\"],\n  [\"80-5=75\", \"72-27=45\"],\n  [\"7+32=39\", \"73-1=72\"],\n  [\"84-12=72\", \"74-14=60\"],\n  [\"33+18=51\", \"95-26=69\"],\n  [\"63-0=63\", \"21+57=78\"],\n  [\"28+54=82\", \"99-55=44\"],\n  [\"22-19=3\", \"49+24=73\"],\n  [\"95-16=79\", \"46+6=52\"],\n  [\"68+1=69\", \"30+19=49\"],\n  [\"46-43=3\", \"29+22=51\"],\n  [\"1+94=95\", \"48+32=80\"],\n  [\"87-6=81\", \"16+41=57\"],\n  [\"80-32=48\", \"49+5=54\"],\n  [\"42+21=63\", \"16+8=24\"],\n  [\"72-29=43\", \"96-8=88\"],\n  [\"81-61=20\", \"56-14=42\"],\n  [\"62+4=66\", \"55+26=81\"],\n  [\"90-76=14\", \"40-3=37\"],\n  [\"21-15=6\", \"67-15=52\"],\n  [\"55+0=55\", \"9+78=87\"],\n  [\"80-11=69\", \"96-36=60\"],\n  [\"95+3=98\", \"73-69=4\"],\n  [\"1+65=66\", \"32+54=86\"],\n  [\"13+10=23\", \"61-59=2\"],\n  [\"1+51=52\", \"36+41=77\"],\n  [\"63-11=52\", \"82-13=69\"],\n  [\"33+10=43\", \"61-36=25\"],\n  [\"73+22=95\", \"22+72=94\"],\n  [\"33-15=18\", \"26+44=70\"],\n  [\"12-2=10\", \"26+8=34\"],\n  [\"49+23=72\", \"67-4=63\"],\n  [\"16-10=6\", \"50-47=3\"],\n  [\"30+49=79\", \"92-9=83\"],\n  [\"24+35=59\", \"72-16=56\"],\n  [\"99-21=78\", \"38+35=73\"],\n  [\"49+12=61\", \"85-32=53\"],\n  [\"4+1=5\", \"37+13=50\"],\n  [\"98-25=73\", \"38+36=74\"],\n  [\"77+17=94\", \"39+32=71\"],\n  [\"39+50=89\", \"47+18=65\"],\n  [\"71+3=74\", \"57-13=44\"],\n  [\"62-29=33\", \"86-40=46\"],\n  [\"90+1=91\", \"11+49=60\"],\n  [\"52+38=90\", \"14+81=95\"],\n  [\"62-42=20\", \"61-31=30\"],\n  [\"98-88=10\", \"48+37=85\"],\n  [\"13-7=6\", \"0+31=31\"],\n  [\"71-5=66\", \"47+31=78\"],\n  [\"37-29=8\", \"54-32=22\"],\n  [\"65+21=86\", \"75-49=26\"],\n  [\"52-10=42\", \"42+19=61\"],\n  [\"99-2=97\", \"87-52=35\"],\n  [\"46-23=23\", \"80-41=39\"],\n  [\"11+23=34\", \"27-16=11\"],\n  [\"41-20=21\", \"87-72=15\"],\n  [\"0+7=7\", \"6-2=4\"],\n  [\"32+10=42\", \"60-37=23\"],\n  [\"17+7=24\", \"48+44=92\"],\n  [\"12+53=65\", \"36+46=82\"],\n  [\"59-46=13\", \"44+43=87\"],\n  [\"42+37=79\", \"8+12=20\"],\n  [\"9+16=25\", \"6+29=35\"],\n  [\"30-23=7\", \"26-10=16\"],\n  [\"25+40=65\", \"84-60=24\"],\n  [\"99-58=41\", \"65-17=48\"],\n  [\"22+31=53\", \"4+85=89\"],\n  [\"90-64=26\", \"84+6=90\"],\n  [\"17+75=92\", \"8+11=19\"],\n  [\"43+0=43\", \"72-65=7\"],\n  [\"94-89=5\", \"49+42=91\"],\n  [\"64-61=3\", \"8+41=49\"],\n  [\"65-33=32\", \"43+52=95\"],\n  [\"18+4=22\", \"8-1=7\"],\n  [\"98-26=72\", \"75-27=48\"],\n  [\"46-9=37\", \"34+40=74\"],\n  [\"96-10=86\", \"14-0=14\"],\n  [\"66-33=33\", \"96-40=56\"],\n  [\"86-64=22\", \"99-95=4\"],\n  [\"92-6=86\", \"48+42=90\"],\n  [\"91-87=4\", \"36+27=63\"],\n  [\"91-62=29\", \"6+27=33\"],\n  [\"49+21=70\", \"87-49=38\"],\n  [\"77-45=32\", \"87+9=96\"],\n  [\"84-62=22\", \"78-45=33\"],\n  [\"80-36=44\", \"31+42=73\"],\n  [\"93-71=22\", \"31+49=80\"],\n  [\"59+17=76\", \"12+21=33\"],\n  [\"60-14=46\", \"0+84=84\"],\n  [\"17-1=16\", \"48+46=94\"],\n  [\"84+14=98\", \"71-43=28\"],\n  [\"97-55=42\", \"80-73=7\"],\n  [\"84-45=39\", \"75-39=36\"],\n  [\"51-4=47\", \"77-43=34\"],\n  [\"67-16=51\", \"48+0=48\"],\n  [\"37-37=0\", \"46-2=44\"],\n  [\"34-23=11\", \"3+88=91\"],\n  [\"73-20=53\", \"90-75=15\"],\n  [\"89+1=90\", \"8+78=86\"],\n  [\"44+55=99\", \"35+19=54\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(`Expected exactly 1 match for \"${oldText}\", found ${results.items.length}`);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n$pairs = @(\n    @(\"2024-09-11 Wednesday\", \"2024-09-12 Thursday\"),\n    @(\"61-55=6\", \"0+11=11\"),\n    @(\"80-5=75\", \"72-27=45\"),\n    @(\"7+32=39\", \"73-1=72\"),\n    @(\"84-12=72\", \"74-14=60\"),\n    @(\"33+18=51\", \"95-26=69\"),\n    @(\"63-0=63\", \"21+57=78\"),\n    @(\"28+54=82\", \"99-55=44\"),\n    @(\"22-19=3\", \"49+24=73\"),\n    @(\"95-16=79\", \"46+6=52\"),\n    @(\"68+1=69\", \"30+19=49\"),\n    @(\"46-43=3\", \"29+22=51\"),\n    @(\"1+94=95\", \"48+32=80\"),\n    @(\"87-6=81\", \"16+41=57\"),\n    @(\"80-32=48\", \"49+5=54\"),\n    @(\"42+21=63\", \"16+8=24\"),\n    @(\"72-29=43\", \"96-8=88\"),\n    @(\"81-61=20\", \"56-14=42\"),\n    @(\"62+4=66\", \"55+26=81\"),\n    @(\"90-76=14\", \"40-3=37\"),\n    @(\"21-15=6\", \"67-15=52\"),\n    @(\"55+0=55\", \"9+78=87\"),\n    @(\"80-11=69\", \"96-36=60\"),\n    @(\"95+3=98\", \"73-69=4\"),\n    @(\"1+65=66\", \"32+54=86\"),\n    @(\"13+10=23\", \"61-59=2\"),\n    @(\"1+51=52\", \"36+41=77\"),\n    @(\"63-11=52\", \"82-13=69\"),\n    @(\"33+10=43\", \"61-36=25\"),\n    @(\"73+22=95\", \"22+72=94\"),\n    @(\"33-15=18\", \"26+44=70\"),\n    @(\"12-2=10\", \"26+8=34\"),\n    @(\"49+23=72\", \"67-4=63\"),\n    @(\"16-10=6\", \"50-47=3\"),\n    @(\"30+49=79\", \"92-9=83\"),\n    @(\"24+35=59\", \"72-16=56\"),\n    @(\"99-21=78\", \"38+35=73\"),\n    @(\"49+12=61\", \"85-32=53\"),\n    @(\"4+1=5\", \"37+13=50\"),\n    @(\"98-25=73\", \"38+36=74\"),\n    @(\"77+17=94\", \"39+32=71\"),\n    @(\"39+50=89\", \"47+18=65\"),\n    @(\"71+3=74\", \"57-13=44\"),\n    @(\"62-29=33\", \"86-40=46\"),\n    @(\"90+1=91\", \"11+49=60\"),\n    @(\"52+38=90\", \"14+81=95\"),\n    @(\"62-42=20\", \"61-31=30\"),\n    @(\"98-88=10\", \"48+37=85\"),\n    @(\"13-7=6\", \"0+31=31\"),\n    @(\"71-5=66\", \"47+31=78\"),\n    @(\"37-29=8\", \"54-32=22\"),\n    @(\"65+21=86\", \"75-49=26\"),\n    @(\"52-10=42\", \"42+19=61\"),\n    @(\"99-2=97\", \"87-52=35\"),\n    @(\"46-23=23\", \"80-41=39\"),\n    @(\"11+23=34\", \"27-16=11\"),\n    @(\"41-20=21\", \"87-72=15\"),\n    @(\"0+7=7\", \"6-2=4\"),\n    @(\"32+10=42\", \"60-37=23\"),\n    @(\"17+7=24\", \"48+44=92\"),\n    @(\"12+53=65\", \"36+46=82\"),\n    @(\"59-46=13\", \"44+43=87\"),\n    @(\"42+37=79\", \"8+12=20\"),\n    @(\"9+16=25\", \"6+29=35\"),\n    @(\"30-23=7\", \"26-10=16\"),\n    @(\"25+40=65\", \"84-60=24\"),\n    @(\"99-58=41\", \"65-17=48\"),\n    @(\"22+31=53\", \"4+85=89\"),\n    @(\"90-64=26\", \"84+6=90\"),\n    @(\"17+75=92\", \"8+11=19\"),\n    @(\"43+0=43\", \"72-65=7\"),\n    @(\"94-89=5\", \"49+42=91\"),\n    @(\"64-61=3\", \"8+41=49\"),\n    @(\"65-33=32\", \"43+52=95\"),\n    @(\"18+4=22\", \"8-1=7\"),\n    @(\"98-26=72\", \"75-27=48\"),\n    @(\"46-9=37\", \"34+40=74\"),\n    @(\"96-10=86\", \"14-0=14\"),\n    @(\"66-33=33\", \"96-40=56\"),\n    @(\"86-64=22\", \"99-95=4\"),\n    @(\"92-6=86\", \"48+42=90\"),\n    @(\"91-87=4\", \"36+27=63\"),\n    @(\"91-62=29\", \"6+27=33\"),\n    @(\"49+21=70\", \"87-49=38\"),\n    @(\"77-45=32\", \"87+9=96\"),\n    @(\"84-62=22\", \"78-45=33\"),\n    @(\"80-36=44\", \"31+42=73\"),\n    @(\"93-71=22\", \"31+49=80\"),\n    @(\"59+17=76\", \"12+21=33\"),\n    @(\"60-14=46\", \"0+84=84\"),\n    @(\"17-1=16\", \"48+46=94\"),\n    @(\"84+14=98\", \"71-43=28\"),\n    @(\"97-55=42\", \"80-73=7\"),\n    @(\"84-45=39\", \"75-39=36\"),\n    @(\"51-4=47\", \"77-43=34\"),\n    @(\"67-16=51\", \"48+0=48\"),\n    @(\"37-37=0\", \"46-2=44\"),\n    @(\"34-23=11\", \"3+88=91\"),\n    @(\"73-20=53\", \"90-75=15\"),\n    @(\"89+1=90\", \"8+78=86\"),\n    @(\"44+55=99\", \"35+19=54\"),\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll)\n    if (-not $found) {\n        throw \"No match found for: $oldText\"\n    }\n}"}
